# Update the embedded build timestamp from "February 03 2026 17.29.55 EST"
# to "February 03 2026 18.05.36 EST" throughout the workbook, as part of
# releasing version 1.0.0 of "Coal Mine Boundaries and Methane Sources".

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

# --- "About" sheet -----------------------------------------------------

$a2 = $wsAbout.Range("A2")
$a2Text = $a2.Value()
$a2.Value = $a2Text.Replace($oldStamp, $newStamp)

$a6 = $wsAbout.Range("A6")
$a6Text = $a6.Value()
$a6.Value = $a6Text.Replace($oldStamp, $newStamp)

# --- "Boundaries and methane sources" sheet -----------------------------

for ($row = 2; $row -le 10; $row++) {
    $cell = $wsData.Cells.Item($row, 19)  # column S ("build_version")
    $cellText = $cell.Value()
    $cell.Value = $cellText.Replace($oldStamp, $newStamp)
}
